# "Community - analyzer updated"
# The analyzer's id_electricity_feed_in value (column K) changes from 1 to 2
# for every data row (rows 2-31), and the sheet's view scrolls/selects a
# different cell (K11 instead of H7, with column F pinned to the left edge).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column K (id_electricity_feed_in) values from 1 to 2 for rows 2-31
$ws.Range("K2:K31").Value = 2

# Update the sheet view: scroll so column F is at the left edge of the
# window, and make K11 the active/selected cell.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K11").Select()
